# unitTest_base_part3.xlsx edit script
# Implements:
#  - [expression] JSON `keys(jsonpath)` / [json] `storeKeys(json,jsonpath,var)` addition
#  - Adjusts the "#system" lookup-table sheet used to drive the function dropdowns:
#      * removes the stray "text" entry from the `target` category list (column A)
#      * inserts a new `storeKeys(json,jsonpath,var)` row into the `json` function list (column M)
#      * the old, 2-row "text" helper column (Y) is dropped, so the `web`/`webalert`/
#        `webcookie`/`ws`/`ws.async`/`xml` columns each slide one letter to the left
#  - Updates the corresponding defined names to match the resized ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# 1) Remove the "text" row from the `target` list in column A (shifts A26:A31 up into A25:A30).
$ws.Range("A25").Delete()

# 2) Insert a new row for `storeKeys(json,jsonpath,var)` into the `json` list in column M,
#    right after storeCount and before storeValue (shifts M16:M17 down into M17:M18).
$ws.Range("M16").Insert()
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# 3) Drop the old 2-entry "text" column (Y); this slides Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD.
$ws.Columns("Y").Delete()

# 4) Re-point the defined names at their new (resized / shifted) ranges.
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

Write-Host "Edit complete"
